# Auto-generated Excel COM-interop script
# Applies cached market-value updates to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 17241924
$ws.Range("I92").Value = 20833800
$ws.Range("J92").Value = 918
$ws.Range("K92").Value = 20833800
$ws.Range("L92").Value = 918
$ws.Range("M92").Value = -20832552
$ws.Range("N92").Value = -3414
# Row 115
$ws.Range("H115").Value = 758.75
$ws.Range("I115").Value = 795
$ws.Range("J115").Value = 650
$ws.Range("K115").Value = 2385
$ws.Range("L115").Value = 1950
$ws.Range("M115").Value = -818
$ws.Range("N115").Value = -5084
# Row 132
$ws.Range("H132").Value = 200189.7
$ws.Range("I132").Value = 4237.977
$ws.Range("J132").Value = 1431886.2
$ws.Range("K132").Value = 12713.931
$ws.Range("L132").Value = 4295658.6
$ws.Range("M132").Value = -10183.931
$ws.Range("N132").Value = -4300718.6
# Row 137
$ws.Range("H137").Value = 4836.1035
$ws.Range("I137").Value = 1022.5
$ws.Range("J137").Value = 8395.467000000001
$ws.Range("K137").Value = 3067.5
$ws.Range("L137").Value = 25186.401
$ws.Range("M137").Value = -517.5
$ws.Range("N137").Value = -30286.401

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4556.0713
$ws.Range("I2").Value = 1031.3478
$ws.Range("J2").Value = 20769.8
$ws.Range("K2").Value = 1031.3478
$ws.Range("L2").Value = 20769.8
$ws.Range("M2").Value = -918.3478
$ws.Range("N2").Value = -20995.8
# Row 32
$ws.Range("H32").Value = 20491.572
$ws.Range("I32").Value = 13587.816
$ws.Range("J32").Value = 34775.207
$ws.Range("K32").Value = 13587.816
$ws.Range("L32").Value = 34775.207
$ws.Range("M32").Value = -13300.816
$ws.Range("N32").Value = -35349.207
# Row 80
$ws.Range("H80").Value = 20857.143
$ws.Range("J80").Value = 20857.143
$ws.Range("L80").Value = 20857.143
$ws.Range("N80").Value = -22853.143
# Row 83
$ws.Range("H83").Value = 20857.143
$ws.Range("J83").Value = 20857.143
$ws.Range("L83").Value = 62571.429
$ws.Range("N83").Value = -72555.429
# Row 97
$ws.Range("H97").Value = 801.71875
$ws.Range("I97").Value = 433.31818
$ws.Range("J97").Value = 1612.2
$ws.Range("K97").Value = 433.31818
$ws.Range("L97").Value = 1612.2
$ws.Range("M97").Value = 62.68182000000002
$ws.Range("N97").Value = -2604.2
# Row 116
$ws.Range("H116").Value = 4556.0713
$ws.Range("I116").Value = 1031.3478
$ws.Range("J116").Value = 20769.8
$ws.Range("K116").Value = 1031.3478
$ws.Range("L116").Value = 20769.8
$ws.Range("M116").Value = 1262.6522
$ws.Range("N116").Value = -25357.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4556.0713
$ws.Range("I3").Value = 1031.3478
$ws.Range("J3").Value = 20769.8
$ws.Range("K3").Value = 1031.3478
$ws.Range("L3").Value = 20769.8
$ws.Range("M3").Value = -917.3478
$ws.Range("N3").Value = -20997.8
# Row 107
$ws.Range("H107").Value = 1025.68
$ws.Range("I107").Value = 840.1429000000001
$ws.Range("J107").Value = 1999.75
$ws.Range("K107").Value = 840.1429000000001
$ws.Range("L107").Value = 1999.75
$ws.Range("M107").Value = 1079.8571
$ws.Range("N107").Value = -5839.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 901.8333
$ws.Range("I16").Value = 555.5
$ws.Range("J16").Value = 1075
$ws.Range("K16").Value = 555.5
$ws.Range("L16").Value = 1075
$ws.Range("M16").Value = -268.5
$ws.Range("N16").Value = -1649
# Row 31
$ws.Range("H31").Value = 18158.922
$ws.Range("I31").Value = 21010.234
$ws.Range("J31").Value = 12565.962
$ws.Range("K31").Value = 21010.234
$ws.Range("L31").Value = 12565.962
$ws.Range("M31").Value = -20715.234
$ws.Range("N31").Value = -13155.962
# Row 34
$ws.Range("H34").Value = 18158.922
$ws.Range("I34").Value = 21010.234
$ws.Range("J34").Value = 12565.962
$ws.Range("K34").Value = 21010.234
$ws.Range("L34").Value = 12565.962
$ws.Range("M34").Value = -20808.234
$ws.Range("N34").Value = -12969.962
# Row 107
$ws.Range("H107").Value = 441.88235
$ws.Range("I107").Value = 319
$ws.Range("J107").Value = 617.4286
$ws.Range("K107").Value = 319
$ws.Range("L107").Value = 617.4286
$ws.Range("M107").Value = 1601
$ws.Range("N107").Value = -4457.4286
# Row 108
$ws.Range("H108").Value = 53190.668
$ws.Range("I108").Value = 20000
$ws.Range("J108").Value = 69786
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 69786
$ws.Range("M108").Value = -16160
$ws.Range("N108").Value = -77466
# Row 113
$ws.Range("H113").Value = 901.8333
$ws.Range("I113").Value = 555.5
$ws.Range("J113").Value = 1075
$ws.Range("K113").Value = 555.5
$ws.Range("L113").Value = 1075
$ws.Range("M113").Value = 1614.5
$ws.Range("N113").Value = -5415
# Row 141
$ws.Range("H141").Value = 43419.453
$ws.Range("I141").Value = 19699
$ws.Range("J141").Value = 46691.242
$ws.Range("K141").Value = 19699
$ws.Range("L141").Value = 46691.242
$ws.Range("M141").Value = -14519
$ws.Range("N141").Value = -57051.242

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 200400.4
$ws.Range("I68").Value = 1000002
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 3000006
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -2999195
$ws.Range("N68").Value = -3122
# Row 71
$ws.Range("H71").Value = 200400.4
$ws.Range("I71").Value = 1000002
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 9000018
$ws.Range("L71").Value = 4500
$ws.Range("M71").Value = -8995962
$ws.Range("N71").Value = -12612
# Row 131
$ws.Range("H131").Value = 115776.63
$ws.Range("I131").Value = 445.3846
$ws.Range("J131").Value = 136037.53
$ws.Range("K131").Value = 1336.1538
$ws.Range("L131").Value = 408112.59
$ws.Range("M131").Value = 3703.8462
$ws.Range("N131").Value = -418192.59

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 6900.7915
$ws.Range("I80").Value = 2700.9092
$ws.Range("J80").Value = 10454.538
$ws.Range("K80").Value = 2700.9092
$ws.Range("L80").Value = 10454.538
$ws.Range("M80").Value = -1702.9092
$ws.Range("N80").Value = -12450.538
# Row 83
$ws.Range("H83").Value = 6900.7915
$ws.Range("I83").Value = 2700.9092
$ws.Range("J83").Value = 10454.538
$ws.Range("K83").Value = 13504.546
$ws.Range("L83").Value = 52272.69
$ws.Range("M83").Value = -8512.546
$ws.Range("N83").Value = -62256.69
# Row 107
$ws.Range("H107").Value = 6436.125
$ws.Range("I107").Value = 256.125
$ws.Range("J107").Value = 12616.125
$ws.Range("K107").Value = 256.125
$ws.Range("L107").Value = 12616.125
$ws.Range("M107").Value = 1663.875
$ws.Range("N107").Value = -16456.125
# Row 113
$ws.Range("H113").Value = 2540
$ws.Range("I113").Value = 3400
$ws.Range("J113").Value = 1966.6666
$ws.Range("K113").Value = 3400
$ws.Range("L113").Value = 1966.6666
$ws.Range("M113").Value = -1230
$ws.Range("N113").Value = -6306.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 685.65216
$ws.Range("I22").Value = 586.8889
$ws.Range("J22").Value = 749.1429000000001
$ws.Range("K22").Value = 586.8889
$ws.Range("L22").Value = 749.1429000000001
$ws.Range("M22").Value = -291.8889
$ws.Range("N22").Value = -1339.1429
# Row 27
$ws.Range("H27").Value = 685.65216
$ws.Range("I27").Value = 586.8889
$ws.Range("J27").Value = 749.1429000000001
$ws.Range("K27").Value = 586.8889
$ws.Range("L27").Value = 749.1429000000001
$ws.Range("M27").Value = -479.8889
$ws.Range("N27").Value = -963.1429000000001
# Row 43
$ws.Range("H43").Value = 29150
$ws.Range("J43").Value = 29150
$ws.Range("L43").Value = 29150
$ws.Range("N43").Value = -29536
# Row 46
$ws.Range("H46").Value = 1403.9231
$ws.Range("I46").Value = 497.77777
$ws.Range("J46").Value = 1883.6471
$ws.Range("K46").Value = 497.77777
$ws.Range("L46").Value = 1883.6471
$ws.Range("M46").Value = -309.77777
$ws.Range("N46").Value = -2259.6471
# Row 55
$ws.Range("H55").Value = 226.7
$ws.Range("I55").Value = 102.30769
$ws.Range("J55").Value = 457.7143
$ws.Range("K55").Value = 102.30769
$ws.Range("L55").Value = 457.7143
$ws.Range("M55").Value = 70.69231000000001
$ws.Range("N55").Value = -803.7143
# Row 82
$ws.Range("H82").Value = 1356.6578
$ws.Range("I82").Value = 1110.174
$ws.Range("J82").Value = 1734.6
$ws.Range("K82").Value = 1110.174
$ws.Range("L82").Value = 1734.6
$ws.Range("M82").Value = -749.174
$ws.Range("N82").Value = -2456.6
# Row 85
$ws.Range("H85").Value = 1356.6578
$ws.Range("I85").Value = 1110.174
$ws.Range("J85").Value = 1734.6
$ws.Range("K85").Value = 1110.174
$ws.Range("L85").Value = 1734.6
$ws.Range("M85").Value = 137.826
$ws.Range("N85").Value = -4230.6
# Row 93
$ws.Range("H93").Value = 1575.6666
$ws.Range("I93").Value = 1220.52
$ws.Range("J93").Value = 2209.8572
$ws.Range("K93").Value = 1220.52
$ws.Range("L93").Value = 2209.8572
$ws.Range("M93").Value = 27.48000000000002
$ws.Range("N93").Value = -4705.8572

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1021.2941
$ws.Range("I100").Value = 562.9091
$ws.Range("J100").Value = 1861.6666
$ws.Range("K100").Value = 1125.8182
$ws.Range("L100").Value = 3723.3332
$ws.Range("M100").Value = -584.8181999999999
$ws.Range("N100").Value = -4805.3332
# Row 107
$ws.Range("H107").Value = 251.61539
$ws.Range("I107").Value = 217.1
$ws.Range("K107").Value = 651.3
$ws.Range("M107").Value = 1268.7

